$d = $word.ActiveDocument
$find = $d.Content.Find
$found = $find.Execute("The context provided does not include specific information about the average sales in Delhi. It mentions cumulative sales figures for Mumbai, Pune, and Bengaluru, as well as overall sales and unit launches in 2024, but it does not break down these figures by city, including Delhi. Therefore, I cannot determine the average sales in Delhi based on the given information.", $true, $false, $false, $false, $false, $true, 1, $false, "The different pH values used by students, as mentioned in the context, are:^l^l- Lemon juice: 2.4, 2.0, 2.2^l- Baking soda (1 Tbsp) in Water (1 cup): 8.4, 8.3, 8.7^l- Orange juice: 3.5, 4.0, 3.4^l- Battery acid: 1.0, 0.7, 0.5^l- Apples: 3.0, 3.2, 3.5^l- Tomatoes: 4.5, 4.2, 4.0^l- Bottled water: 6.7, 7.0, 7.2^l- Milk of magnesia: 10.5, 10.3, 10.6^l- Liquid hand soap: 9.0, 10.0, 9.5^l- Vinegar: 2.2, 2.9, 3.0^l- Household bleach: 12.5, 12.5, 12.7^l- Milk: 6.6, 6.5, 6.4^l- Household ammonia: 11.5, 11.0, 11.5^l- Lye: 13.0, 13.5, 13.4^l- Sodium hydroxide: 14.0, 14.0, 13.9^l- Anti-freeze: 10.1, 10.9, 9.7^l- Windex: 9.9, 10.2, 9.5^l- Liquid detergent: 10.5, 10.0, 10.3^l- Cola: 3.0, 2.5, 3.2", 2)
Write-Output $found
